# Update the roster table so that rows 3 and 11-16 reflect the new
# player / position / team assignments (players shuffled to new rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Shaedon Sharpe (SG,SF / Portland Trail Blazers) -> Mikal Bridges (SG,SF,PF / New York Knicks)
$ws.Range("A3").Value = "Mikal Bridges"
$ws.Range("B3").Value = "SG,SF,PF"
$ws.Range("C3").Value = "New York Knicks"

# Row 11: Clint Capela (C / Atlanta Hawks) -> Brook Lopez (C / Milwaukee Bucks)
$ws.Range("A11").Value = "Brook Lopez"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Milwaukee Bucks"

# Row 12: Brook Lopez (C / Milwaukee Bucks) -> Zach Edey (C / Memphis Grizzlies)
$ws.Range("A12").Value = "Zach Edey"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "Memphis Grizzlies"

# Row 13: Josh Giddey (PG,SG,SF / Chicago Bulls) -> Scottie Barnes (PG,SG,SF,PF / Toronto Raptors)
$ws.Range("A13").Value = "Scottie Barnes"
$ws.Range("B13").Value = "PG,SG,SF,PF"
$ws.Range("C13").Value = "Toronto Raptors"

# Row 14: Scottie Barnes (PG,SG,SF,PF / Toronto Raptors) -> Josh Giddey (PG,SG,SF / Chicago Bulls)
$ws.Range("A14").Value = "Josh Giddey"
$ws.Range("B14").Value = "PG,SG,SF"
$ws.Range("C14").Value = "Chicago Bulls"

# Row 15: Zach Edey (C / Memphis Grizzlies) -> Clint Capela (C / Atlanta Hawks)
$ws.Range("A15").Value = "Clint Capela"
$ws.Range("B15").Value = "C"
$ws.Range("C15").Value = "Atlanta Hawks"

# Row 16: Mikal Bridges (SG,SF,PF / New York Knicks) -> Shaedon Sharpe (SG,SF / Portland Trail Blazers)
$ws.Range("A16").Value = "Shaedon Sharpe"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "Portland Trail Blazers"
